# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the crypto ticker rows per the commit's refreshed data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellRefs = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "E25",
    "E26",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "E41",
    "D42",
    "E42",
    "E43",
    "D44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D50",
    "E50",
    "D51",
    "E51"
)

$newValues = @(
    '308.38',
    '0.58%',
    '40.98',
    '1.22%',
    '1.53%',
    '0.07621',
    '0.32%',
    '4.278',
    '0.86%',
    '0.9096',
    '0.56%',
    '0.1228',
    '25.52%',
    '0.1809',
    '2.85%',
    '0.08989',
    '-2.13%',
    '0.04273',
    '-2.29%',
    '0.1045',
    '-0.79%',
    '0.001256',
    '1.98%',
    '0.005915',
    '1.63%',
    '3.351',
    '-0.66%',
    '6.936',
    '1.65%',
    '0.1394',
    '3.21%',
    '0.2706',
    '-4.97%',
    '0.04045',
    '-2.84%',
    '0.001272',
    '4.55%',
    '0.004124',
    '1.42%',
    '-2.31%',
    '24.49%',
    '0.02412',
    '0.26%',
    '0.05221',
    '1.88%',
    '0.007842',
    '-0.15%',
    '-0.13%',
    '0.006803',
    '-3.71%',
    '-0.93%',
    '0.008086',
    '0.3063',
    '-7.66%',
    '0.00006899',
    '7.00%',
    '0.00000000751',
    '-0.07%',
    '0.1080',
    '1,605.64%',
    '0.00002104',
    '-0.07%',
    '0.0002003',
    '-0.07%'
)

for ($i = 0; $i -lt $cellRefs.Length; $i++) {
    $rng = $ws.Range($cellRefs[$i])
    # Force text interpretation so values like "308.38" or "0.58%"
    # are stored verbatim as strings, matching the source data feed
    # instead of being auto-parsed into numbers/percentages.
    $rng.NumberFormat = "@"
    $rng.Value = $newValues[$i]
    $rng.Style = "Normal"
}
